$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 82 continues the existing series: copy row 81's formatting (date
# number format / font / border on column A) down to the new row, then
# fill in the new date serial and value.
$ws.Range("A81:B81").Copy()
$ws.Range("A82:B82").PasteSpecial(-4122)

$ws.Range("A82").Value = 45884
$ws.Range("B82").Value = 105.26
